$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.019164323806763
$ws.Range("B1").Value = 3.282618045806885
$ws.Range("C1").Value = 3.772308826446533
$ws.Range("D1").Value = 2.006429672241211
$ws.Range("E1").Value = 1.183880805969238
